# Update cryptos list prices (column D) and 1h volume percentages (column E).
# Values that look like plain numbers are prefixed with a leading apostrophe
# so Excel stores them as text (matching the original inlineStr string type)
# instead of silently converting them to numeric cells; the style is reset
# back to Normal afterwards so no stray "quote prefix" formatting sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.762.36"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.595.16"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("D5").Value = "'209.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("E6").Value = "  -1.86%  "

$ws.Range("E8").Value = "  -2.36%  "

$ws.Range("D11").Value = "'0.0869"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "1.821.69"
$ws.Range("E12").Value = "  -1.48%  "

$ws.Range("D13").Value = "1.587.96"
$ws.Range("E13").Value = "  -1.94%  "

$ws.Range("E14").Value = "  -2.33%  "

$ws.Range("E15").Value = "  -3.35%  "

$ws.Range("D16").Value = "27.760.52"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").Value = "'63.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").Value = "'219.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.89%  "

$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.55%  "

$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("E22").Value = "  -3.10%  "

$ws.Range("D23").Value = "'9.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.20%  "

$ws.Range("E24").Value = "  -3.70%  "

$ws.Range("D25").Value = "'154.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").Value = "'7.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.75%  "

$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").Value = "'15.17"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").Value = "'0.106"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.21%  "

$ws.Range("E30").Value = "  -0.74%  "

$ws.Range("D31").Value = "'0.0473"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("E32").Value = "  -3.83%  "

$ws.Range("D33").Value = "1.379.64"
$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").Value = "'2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("D36").Value = "'0.978"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("E40").Value = "  -1.58%  "

$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").Value = "'0.978"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.87%  "

$ws.Range("E43").Value = "  -0.97%  "

$ws.Range("E44").Value = "  +2.69%  "

$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("E46").Value = "  -2.11%  "

$ws.Range("D47").Value = "1.732.75"
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("D48").Value = "'86.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.62%  "

$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("D50").Value = "'0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.72%  "

$ws.Range("E51").Value = "  -1.10%  "
